$d = $word.ActiveDocument

# This document is a flat list of one-line "TODO"/code-snippet paragraphs.
# Several of them were originally typed with spell-check markers
# (<w:proofErr>) splitting a single logical line across multiple <w:r>
# runs. The edit collapses each of those back into one clean run, and
# appends a new " (x)" run (as a separate <w:r>, matching how Word marks
# freshly-typed text distinctly from pre-existing text) to six of the
# TODO-list lines. The trailing "_GoBack" bookmark - which Word drops at
# the location of the most recent edit - moves from the last paragraph
# ("BubbleTon") to the second-to-last one ("adjustSize") accordingly.

function Rebuild-Paragraph($index, $text) {
    # Replaces paragraph $index's entire content with a single clean run
    # containing $text, removing any proofErr marks / run splits along the
    # way (delete the whole paragraph range, then retype it from scratch).
    $p = $d.Paragraphs.Item($index)
    $full = $p.Range
    $full.Delete()
    $p2 = $d.Paragraphs.Item($index)
    $p2.Range.InsertBefore($text + "`r")
}

function Append-NewRun($index, $text) {
    # Appends $text as a brand-new, separate run at the end of paragraph
    # $index (before the paragraph mark). A plain InsertAfter would just
    # extend the previous run's text (same formatting => coalesced on
    # save), so we break the seam with a transient bookmark: inserting
    # and immediately deleting a zero-length bookmark at the junction
    # forces the engine to keep the two text spans in distinct <w:r>
    # elements without leaving any bookmark/rPr residue behind.
    $p = $d.Paragraphs.Item($index)
    $insertPos = $p.Range.End - 1
    $r = $d.Range($insertPos, $insertPos)
    $r.InsertAfter($text)
    $seam = $d.Range($insertPos, $insertPos)
    $d.Bookmarks.Add("__seam__", $seam)
    $d.Bookmarks("__seam__").Delete()
}

function Set-GoBackBookmarkAfterParagraph($index) {
    # Places a collapsed "_GoBack" bookmark right at the end of paragraph
    # $index (before its paragraph mark), removing it from wherever it
    # was before (_GoBack is a singleton bookmark).
    $p = $d.Paragraphs.Item($index)
    $pos = $p.Range.End - 1
    $placeholder = $d.Range($pos, $pos)
    $placeholder.InsertAfter("Z")
    $bmRange = $d.Range($pos, $pos + 1)
    $d.Bookmarks.Add("_GoBack", $bmRange)
    $delRange = $d.Range($pos, $pos + 1)
    $delRange.Delete()
}

# --- Code snippet paragraphs: collapse proofErr-split runs into single runs ---
Rebuild-Paragraph 5 "getPlatStyleAt:"
Rebuild-Paragraph 6 "if (this.mg[i].x < pos && this.mg[i].x + this.bg[0].image.width >= pos) {"
Rebuild-Paragraph 8 "updatePlatforms:"
Rebuild-Paragraph 9 "if (this.platforms[i].start < -this.plat.segmentSize * (game.tiles_max+3) ) {"
Rebuild-Paragraph 12 "counter -= 0.5 * Math.max(1,this.distance/700);"

# --- TODO list paragraphs: clean up + append " (x)" as a distinct run ---
Rebuild-Paragraph 16 "mehrere Canvas"
Append-NewRun 16 " (x)"

Append-NewRun 17 " (x)"

Append-NewRun 18 " (x)"

Rebuild-Paragraph 19 "richtig adden"
Append-NewRun 19 " (x)"

Rebuild-Paragraph 20 "adjustSize"
Append-NewRun 20 " (x)"

Append-NewRun 21 " (x)"

# --- Move the _GoBack bookmark from the BubbleTon paragraph to the end of
#     the adjustSize paragraph ---
Set-GoBackBookmarkAfterParagraph 20
